$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5933
$ws.Range("E2").Value = 99
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 31
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 48
$ws.Range("J2").Value = -41
$ws.Range("K2").Value = 4463
$ws.Range("L2").Value = 2533
$ws.Range("M2").Value = 1930
$ws.Range("N2").Value = 1882
$ws.Range("O2").Value = 48
$ws.Range("P2").Value = 76
$ws.Range("Q2").Value = 244
$ws.Range("R2").Value = -269
$ws.Range("S2").Value = 44
$ws.Range("T2").Value = 263
$ws.Range("U2").Value = -19
$ws.Range("V2").Value = 1271
$ws.Range("W2").Value = 1.67
$ws.Range("X2").Value = 0.13
$ws.Range("Y2").Value = 2.59
$ws.Range("Z2").Value = 0.17
$ws.Range("AA2").Value = 131.22
$ws.Range("AB2").Value = 2393.68
$ws.Range("AC2").Value = 320
$ws.Range("AD2").Value = 15.32
$ws.Range("AE2").Value = 12963
$ws.Range("AF2").Value = 0.38
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 2.04
$ws.Range("AI2").Value = 29.95
$ws.Range("AJ2").Value = 15160128

# Row 3
$ws.Range("D3").Value = 6595
$ws.Range("E3").Value = 182
$ws.Range("F3").Value = 182
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 87
$ws.Range("I3").Value = 118
$ws.Range("J3").Value = -31
$ws.Range("K3").Value = 4721
$ws.Range("L3").Value = 2695
$ws.Range("M3").Value = 2026
$ws.Range("N3").Value = 2011
$ws.Range("O3").Value = 15
$ws.Range("P3").Value = 76
$ws.Range("Q3").Value = 419
$ws.Range("R3").Value = -288
$ws.Range("S3").Value = -78
$ws.Range("T3").Value = 262
$ws.Range("U3").Value = 157
$ws.Range("V3").Value = 1266
$ws.Range("W3").Value = 2.76
$ws.Range("X3").Value = 1.32
$ws.Range("Y3").Value = 6.08
$ws.Range("Z3").Value = 1.9
$ws.Range("AA3").Value = 132.98
$ws.Range("AB3").Value = 2543.5
$ws.Range("AC3").Value = 781
$ws.Range("AD3").Value = 7.2
$ws.Range("AE3").Value = 13853
$ws.Range("AF3").Value = 0.41
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 2.67
$ws.Range("AI3").Value = 18.4
$ws.Range("AJ3").Value = 15160128

# Row 4
$ws.Range("D4").Value = 6240
$ws.Range("E4").Value = 184
$ws.Range("F4").Value = 184
$ws.Range("G4").Value = 148
$ws.Range("H4").Value = 115
$ws.Range("I4").Value = 102
$ws.Range("J4").Value = 14
$ws.Range("K4").Value = 4803
$ws.Range("L4").Value = 2725
$ws.Range("M4").Value = 2078
$ws.Range("N4").Value = 1972
$ws.Range("O4").Value = 107
$ws.Range("P4").Value = 76
$ws.Range("Q4").Value = 261
$ws.Range("R4").Value = -243
$ws.Range("S4").Value = 70
$ws.Range("T4").Value = 214
$ws.Range("U4").Value = 47
$ws.Range("V4").Value = 1393
$ws.Range("W4").Value = 2.94
$ws.Range("X4").Value = 1.85
$ws.Range("Y4").Value = 5.1
$ws.Range("Z4").Value = 2.42
$ws.Range("AA4").Value = 131.12
$ws.Range("AB4").Value = 2620.24
$ws.Range("AC4").Value = 670
$ws.Range("AD4").Value = 7.67
$ws.Range("AE4").Value = 13578
$ws.Range("AF4").Value = 0.38
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 1.95
$ws.Range("AI4").Value = 14.29
$ws.Range("AJ4").Value = 15160128

# Row 5
$ws.Range("D5").Value = 5409
$ws.Range("E5").Value = 50
$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 38
$ws.Range("H5").Value = -9
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = -9
$ws.Range("K5").Value = 4569
$ws.Range("L5").Value = 2556
$ws.Range("M5").Value = 2013
$ws.Range("N5").Value = 1921
$ws.Range("O5").Value = 92
$ws.Range("P5").Value = 76
$ws.Range("Q5").Value = 468
$ws.Range("R5").Value = -328
$ws.Range("S5").Value = -13
$ws.Range("T5").Value = 257
$ws.Range("U5").Value = 210
$ws.Range("V5").Value = 1430
$ws.Range("W5").Value = 0.92
$ws.Range("X5").Value = -0.17
$ws.Range("Y5").Value = -0.01
$ws.Range("Z5").Value = -0.2
$ws.Range("AA5").Value = 126.97
$ws.Range("AB5").Value = 2614.64
$ws.Range("AC5").Value = -1
$ws.Range("AD5").Value = -5723.17
$ws.Range("AE5").Value = 13230
$ws.Range("AF5").Value = 0.43
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 2.64
$ws.Range("AI5").Value = -14475.46
$ws.Range("AJ5").Value = 15160128

# Row 6
$ws.Range("D6").Value = 5041
$ws.Range("E6").Value = 111
$ws.Range("F6").Value = 111
$ws.Range("G6").Value = 48
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 30
$ws.Range("K6").Value = 4554
$ws.Range("L6").Value = 2576
$ws.Range("M6").Value = 1978
$ws.Range("N6").Value = 1903
$ws.Range("P6").Value = 76
$ws.Range("Q6").Value = 79
$ws.Range("R6").Value = -317
$ws.Range("S6").Value = 63
$ws.Range("T6").Value = 351
$ws.Range("U6").Value = -272
$ws.Range("V6").Value = 1563
$ws.Range("W6").Value = 2.2
$ws.Range("X6").Value = 0.19
$ws.Range("Y6").Value = 1.56
$ws.Range("Z6").Value = 0.22
$ws.Range("AA6").Value = 130.2
$ws.Range("AB6").Value = 2588.5
$ws.Range("AC6").Value = 197
$ws.Range("AD6").Value = 34.23
$ws.Range("AE6").Value = 13105
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 1.49
$ws.Range("AI6").Value = 48.71
$ws.Range("AJ6").Value = 15160128

# Clear rows 7-9 data columns D:AJ, leaving A-C intact
$ws.Range("D7:AJ9").ClearContents()
